# Profile onboarding script implementation:
# Append new test case "Profile59" (OPQA-2108) as row 60, matching the
# plain (non-wrapped) formatting already used elsewhere in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the new row's cell formatting from existing "plain" styled cells
# (A46 = label style, E46 = plain value style) so the new row blends in
# with the rest of the test-case table instead of picking up ad-hoc styles.
$ws.Range("A46").Copy()
$ws.Range("A60").PasteSpecial(-4122)
$ws.Range("E46").Copy()
$ws.Range("B60:E60").PasteSpecial(-4122)

# New test case content
$ws.Cells.Item(60, 1).Value = "Profile59"
$ws.Cells.Item(60, 2).Value = "OPQA-2108"
$ws.Cells.Item(60, 3).Value = "Verify that user has the ability to add and update the following information from the profile modal:Title/Role,Institution,Country,Skills and Interests (Topics)"
$ws.Cells.Item(60, 4).Value = "Y"

# Reflect the author's final cursor position on the new row
$ws.Range("D60").Select()
